$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename metrics: mean_ir -> mfrd, cv_ir -> afrd
$ws.Range("C5").Value = "mfrd"
$ws.Range("C6").Value = "afrd"

# Update Mean/Std values for the renamed metrics
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.60398
$ws.Range("E6").Value = 0.06195
